# "Changed S2F from 1m avg to 12m avg"
#
# The author duplicated the "PreHalving" sheet (Excel inserts the copy
# immediately after the original, named "PreHalving (2)"), renamed the
# copy to "PreHalving (SFW)", and then updated the regression inputs
# (intercept / SF-coefficient, cells B1/B2) on both the original
# "PreHalving" sheet and the new "PreHalving (SFW)" sheet to reflect the
# refreshed stock-to-flow regression (12-month average instead of the
# previous 1-month average). A couple of leftover cell selections
# (navigation, not data) also moved on two other sheets.

$wb = $excel.ActiveWorkbook

# 1. Duplicate "PreHalving" -> creates "PreHalving (2)" right after it,
#    then rename it to match the author's new sheet name.
$preHalving = $wb.Worksheets.Item("PreHalving")
$preHalving.Copy($null, $preHalving)
$sfw = $wb.Worksheets.Item(2)
$sfw.Name = "PreHalving (SFW)"

# 2. Refresh the regression inputs on the original "PreHalving" sheet.
$preHalving.Range("B1").Value = 15.37
$preHalving.Range("B2").Value = 3.78

# 3. Refresh the regression inputs on the new "PreHalving (SFW)" sheet
#    (its own, slightly different, re-fit coefficients).
$sfw.Range("B1").Value = 15.38
$sfw.Range("B2").Value = 3.79

# 4. Leftover navigation: the selection on "AllData" moved to G14.
$allData = $wb.Worksheets.Item("AllData")
$allData.Range("G14").Select()

# 5. Leave "PreHalving (SFW)" as the active sheet/tab, selection on G17
#    (matches where the author left the cursor after editing it).
$sfw.Activate()
$sfw.Range("G17").Select()
